$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "In Translation"
#    This string is shared by every cell that shows that status:
#      - Status column (column C) on the "zh-cn" and "de-de" tables,
#        rows 2 and 3
#      - The per-locale summary columns ("zh-cn"/"de-de", i.e.
#        columns E/F) on the "Overview" sheet, rows 2 and 3
#    All of these need to be updated together since they originally
#    point at the same shared string.
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# -----------------------------------------------------------------
# 2) Column width changes (17.2159881591797 -> 13.4101848602295
#    stored/raw XML width). Excel's ColumnWidth property is expressed
#    in "characters" and gets snapped to whole pixels on write
#    (stored width = floor(ColumnWidth*6 + 5.5) / 6), so a
#    ColumnWidth of 12.5 is the value that lands on the pixel bucket
#    nearest the target stored width.
# -----------------------------------------------------------------
$newColWidth = 12.5

$wsOverview.Range("E1").ColumnWidth = $newColWidth   # zh-cn column
$wsOverview.Range("F1").ColumnWidth = $newColWidth   # de-de column

$wsZhCn.Range("C1").ColumnWidth = $newColWidth       # Status column
$wsDeDe.Range("C1").ColumnWidth = $newColWidth       # Status column
